$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44250
$ws.Range("J2").Value = 38
$ws.Range("K2").Value = 18000
$ws.Range("L2").Value = 18000
$ws.Range("M2").Value = 18000
$ws.Range("N2").Value = '$/malla 25 kilos'
$ws.Range("O2").Value = 'Provincia de Talca'
$ws.Range("P2").Value = 720

# Row 3
$ws.Range("D3").Value = 44253
$ws.Range("K3").Value = 18000
$ws.Range("L3").Value = 18000
$ws.Range("M3").Value = 18000
$ws.Range("N3").Value = '$/saco 25 kilos'
$ws.Range("O3").Value = 'Provincia de Talca'
$ws.Range("P3").Value = 720

# Row 5
$ws.Range("D5").Value = 44452
$ws.Range("J5").Value = 70
$ws.Range("K5").Value = 31000
$ws.Range("L5").Value = 32000
$ws.Range("M5").Value = 31500
$ws.Range("N5").Value = '$/malla 25 kilos'
$ws.Range("O5").Value = 'Provincia de Limarí'
$ws.Range("P5").Value = 1260

# Row 6
$ws.Range("D6").Value = 44475
$ws.Range("J6").Value = 73
$ws.Range("K6").Value = 25000
$ws.Range("L6").Value = 26000
$ws.Range("M6").Value = 25479
$ws.Range("O6").Value = 'Provincia de Limarí'
$ws.Range("P6").Value = 1019

# Row 7
$ws.Range("D7").Value = 44410
$ws.Range("J7").Value = 35
$ws.Range("K7").Value = 34000
$ws.Range("L7").Value = 34000
$ws.Range("M7").Value = 34000
$ws.Range("N7").Value = '$/malla 25 kilos'
$ws.Range("P7").Value = 1360

# Row 8
$ws.Range("D8").Value = 44365
$ws.Range("J8").Value = 70
$ws.Range("K8").Value = 22000
$ws.Range("L8").Value = 23000
$ws.Range("M8").Value = 22500
$ws.Range("P8").Value = 900

# Row 9
$ws.Range("D9").Value = 44448
$ws.Range("J9").Value = 45
$ws.Range("K9").Value = 32000
$ws.Range("L9").Value = 32000
$ws.Range("M9").Value = 32000
$ws.Range("P9").Value = 1280

# Row 10
$ws.Range("D10").Value = 44252
$ws.Range("J10").Value = 40
$ws.Range("K10").Value = 18000
$ws.Range("L10").Value = 19000
$ws.Range("M10").Value = 18625
$ws.Range("N10").Value = '$/malla 25 kilos'
$ws.Range("O10").Value = 'Provincia de Talca'
$ws.Range("P10").Value = 745

# Row 11
$ws.Range("D11").Value = 44165
$ws.Range("J11").Value = 45
$ws.Range("K11").Value = 22000
$ws.Range("L11").Value = 22000
$ws.Range("M11").Value = 22000
$ws.Range("N11").Value = '$/saco 25 kilos'
$ws.Range("O11").Value = 'Provincia de Quillota'
$ws.Range("P11").Value = 880

# Row 12
$ws.Range("D12").Value = 44483
$ws.Range("J12").Value = 55
$ws.Range("K12").Value = 29000
$ws.Range("L12").Value = 30000
$ws.Range("M12").Value = 29455
$ws.Range("P12").Value = 1178

# Row 13
$ws.Range("D13").Value = 44468
$ws.Range("J13").Value = 65
$ws.Range("L13").Value = 25000
$ws.Range("M13").Value = 24538
$ws.Range("P13").Value = 982

# Row 14
$ws.Range("D14").Value = 44161
$ws.Range("J14").Value = 35
$ws.Range("K14").Value = 21000
$ws.Range("L14").Value = 21000
$ws.Range("M14").Value = 21000
$ws.Range("N14").Value = '$/saco 25 kilos'
$ws.Range("O14").Value = 'Provincia de Quillota'
$ws.Range("P14").Value = 840

# Row 15
$ws.Range("D15").Value = 44162
$ws.Range("J15").Value = 35
$ws.Range("K15").Value = 17000
$ws.Range("L15").Value = 17000
$ws.Range("M15").Value = 17000
$ws.Range("O15").Value = 'Provincia de Quillota'
$ws.Range("P15").Value = 680

# Row 16
$ws.Range("D16").Value = 44159
$ws.Range("J16").Value = 35
$ws.Range("K16").Value = 22000
$ws.Range("L16").Value = 22000
$ws.Range("M16").Value = 22000
$ws.Range("N16").Value = '$/malla 25 kilos'
$ws.Range("O16").Value = 'Provincia de Quillota'
$ws.Range("P16").Value = 880

# Row 17
$ws.Range("D17").Value = 44412
$ws.Range("K17").Value = 24000
$ws.Range("L17").Value = 24000
$ws.Range("M17").Value = 24000
$ws.Range("P17").Value = 960

# Row 18
$ws.Range("D18").Value = 44160
$ws.Range("J18").Value = 35
$ws.Range("K18").Value = 21000
$ws.Range("L18").Value = 21000
$ws.Range("M18").Value = 21000
$ws.Range("N18").Value = '$/saco 25 kilos'
$ws.Range("O18").Value = 'Provincia de Quillota'
$ws.Range("P18").Value = 840

# Row 20
$ws.Range("D20").Value = 44399
$ws.Range("J20").Value = 38
$ws.Range("K20").Value = 33000
$ws.Range("L20").Value = 33000
$ws.Range("M20").Value = 33000
$ws.Range("P20").Value = 1320

# Row 21
$ws.Range("D21").Value = 44376
$ws.Range("J21").Value = 38
$ws.Range("K21").Value = 27000
$ws.Range("L21").Value = 27000
$ws.Range("M21").Value = 27000
$ws.Range("O21").Value = 'Provincia de Limarí'
$ws.Range("P21").Value = 1080

# Row 22
$ws.Range("D22").Value = 44453
$ws.Range("J22").Value = 73
$ws.Range("K22").Value = 21000
$ws.Range("M22").Value = 21521
$ws.Range("O22").Value = 'Provincia de Limarí'
$ws.Range("P22").Value = 861

# Row 23
$ws.Range("D23").Value = 44476
$ws.Range("K23").Value = 23000
$ws.Range("L23").Value = 24000
$ws.Range("M23").Value = 23521
$ws.Range("P23").Value = 941

# Row 24
$ws.Range("D24").Value = 44372
$ws.Range("J24").Value = 50
$ws.Range("K24").Value = 33000
$ws.Range("L24").Value = 34000
$ws.Range("M24").Value = 33500
$ws.Range("N24").Value = '$/saco 25 kilos'
$ws.Range("O24").Value = 'Provincia de Limarí'
$ws.Range("P24").Value = 1340

# Row 25
$ws.Range("D25").Value = 44469
$ws.Range("J25").Value = 73
$ws.Range("K25").Value = 28000
$ws.Range("L25").Value = 29000
$ws.Range("M25").Value = 28521
$ws.Range("N25").Value = '$/malla 25 kilos'
$ws.Range("P25").Value = 1141

# Row 26
$ws.Range("D26").Value = 44481
$ws.Range("J26").Value = 63
$ws.Range("K26").Value = 22000
$ws.Range("L26").Value = 23000
$ws.Range("M26").Value = 22476
$ws.Range("N26").Value = '$/saco 25 kilos'
$ws.Range("P26").Value = 899

# Row 27
$ws.Range("D27").Value = 44487
$ws.Range("J27").Value = 73
$ws.Range("K27").Value = 20000
$ws.Range("L27").Value = 21000
$ws.Range("M27").Value = 20521
$ws.Range("P27").Value = 821

# Row 28
$ws.Range("D28").Value = 44343
$ws.Range("J28").Value = 40
$ws.Range("K28").Value = 28000
$ws.Range("L28").Value = 28000
$ws.Range("M28").Value = 28000
$ws.Range("P28").Value = 1120

# Row 29
$ws.Range("D29").Value = 44370
$ws.Range("J29").Value = 45
$ws.Range("K29").Value = 32000
$ws.Range("L29").Value = 32000
$ws.Range("M29").Value = 32000
$ws.Range("O29").Value = 'Provincia de Limarí'
$ws.Range("P29").Value = 1280

# Row 30
$ws.Range("D30").Value = 44473
$ws.Range("J30").Value = 85
$ws.Range("K30").Value = 35000
$ws.Range("L30").Value = 36000
$ws.Range("M30").Value = 35471
$ws.Range("N30").Value = '$/malla 25 kilos'
$ws.Range("P30").Value = 1419

# Row 31
$ws.Range("D31").Value = 44484
$ws.Range("J31").Value = 71
$ws.Range("K31").Value = 29000
$ws.Range("L31").Value = 30000
$ws.Range("M31").Value = 29507
$ws.Range("N31").Value = '$/saco 25 kilos'
$ws.Range("P31").Value = 1180

# Row 32
$ws.Range("D32").Value = 44181
$ws.Range("K32").Value = 26000
$ws.Range("L32").Value = 26000
$ws.Range("M32").Value = 26000
$ws.Range("O32").Value = 'Región Metropolitana'
$ws.Range("P32").Value = 1040
